$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.6478768977981701
$ws.Range("C2").Value = 0.06827765426241683
$ws.Range("D2").Value = 0.07804545101215865
$ws.Range("E2").Value = 0.08603117381810321
$ws.Range("G2").Value = 0.002476208907534175
$ws.Range("I2").Value = 0.9756134643876067
$ws.Range("K2").Value = 0.4213791976436596
$ws.Range("L2").Value = 0.2164816056303493
$ws.Range("N2").Value = 1.860019786817835
$ws.Range("O2").Value = 4.013973054213807

$ws.Range("B3").Value = 0.606460114415313
$ws.Range("C3").Value = 0.06563215932921906
$ws.Range("D3").Value = 0.07086771833321848
$ws.Range("E3").Value = 0.0854988897022082
$ws.Range("G3").Value = 0.002478875706680454
$ws.Range("I3").Value = 0.9804678646478528
$ws.Range("K3").Value = 0.382139645420466
$ws.Range("L3").Value = 0.2091502386510058
$ws.Range("N3").Value = 1.878671632541064
$ws.Range("O3").Value = 4.02390316916339

$ws.Range("B4").Value = 0.5812682110513947
$ws.Range("C4").Value = 0.06398635146023679
$ws.Range("D4").Value = 0.06649461010776747
$ws.Range("E4").Value = 0.08521366835346633
$ws.Range("G4").Value = 0.002480601686933281
$ws.Range("I4").Value = 0.9838942384609481
$ws.Range("K4").Value = 0.3581438944201807
$ws.Range("L4").Value = 0.2047661608942093
$ws.Range("N4").Value = 1.890709128914217
$ws.Range("O4").Value = 4.031994060331755

$ws.Range("B5").Value = 0.5710627906438503
$ws.Range("C5").Value = 0.06331028040702336
$ws.Range("D5").Value = 0.06472110486009797
$ws.Range("E5").Value = 0.08510791786218519
$ws.Range("G5").Value = 0.00248132737211903
$ws.Range("I5").Value = 0.9854026405627891
$ws.Range("K5").Value = 0.3483904007983654
$ws.Range("L5").Value = 0.2030092058002992
$ws.Range("N5").Value = 1.895761646526793
$ws.Range("O5").Value = 4.035792643138137

$ws.Range("B6").Value = 0.5693718590270009
$ws.Range("C6").Value = 0.06319769390435681
$ws.Range("D6").Value = 0.06442713368186048
$ws.Range("E6").Value = 0.08509099154260724
$ws.Range("G6").Value = 0.002481449222257914
$ws.Range("I6").Value = 0.9856598829322607
$ws.Range("K6").Value = 0.3467723621184291
$ws.Range("L6").Value = 0.202719254051857
$ws.Range("N6").Value = 1.896609499873932
$ws.Range("O6").Value = 4.036453686963512

$ws.Range("B7").Value = 0.5811303314931138
$ws.Range("C7").Value = 0.06397725554893441
$ws.Range("D7").Value = 0.06647065727490542
$ws.Range("E7").Value = 0.08521219970994132
$ws.Range("G7").Value = 0.002480611383315496
$ws.Range("I7").Value = 0.9839141272508698
$ws.Range("K7").Value = 0.3580122536903616
$ws.Range("L7").Value = 0.2047423460746103
$ws.Range("N7").Value = 1.890776673183728
$ws.Range("O7").Value = 4.03204325870874

$ws.Range("B8").Value = 0.6335473020588438
$ws.Range("C8").Value = 0.06736994531265594
$ws.Range("D8").Value = 0.07556350584270888
$ws.Range("E8").Value = 0.08583901795745419
$ws.Range("G8").Value = 0.002477110080386039
$ws.Range("I8").Value = 0.977194768000146
$ws.Range("K8").Value = 0.407829425606252
$ws.Range("L8").Value = 0.2139294197674104
$ws.Range("N8").Value = 1.866329483185154
$ws.Range("O8").Value = 4.016983223987552

$ws.Range("B9").Value = 0.7382060200598346
$ws.Range("C9").Value = 0.07385270347283779
$ws.Range("D9").Value = 0.09366522172192049
$ws.Range("E9").Value = 0.08739767858406111
$ws.Range("G9").Value = 0.002470943657743601
$ws.Range("I9").Value = 0.9675537695862744
$ws.Range("K9").Value = 0.5062793425250334
$ws.Range("L9").Value = 0.2328750205098942
$ws.Range("N9").Value = 1.82303153016994
$ws.Range("O9").Value = 4.003268078547791

$ws.Range("B10").Value = 0.8162188877691676
$ws.Range("C10").Value = 0.07851218244519487
$ws.Range("D10").Value = 0.1071316425505557
$ws.Range("E10").Value = 0.08874307218259503
$ws.Range("G10").Value = 0.002466835463797399
$ws.Range("I10").Value = 0.9626250915532069
$ws.Range("K10").Value = 0.5790602153761313
$ws.Range("L10").Value = 0.247360439543769
$ws.Range("N10").Value = 1.794049498310022
$ws.Range("O10").Value = 4.002837122963768

$ws.Range("B11").Value = 0.8519485623899357
$ws.Range("C11").Value = 0.08060955844545958
$ws.Range("D11").Value = 0.1132946669962394
$ws.Range("E11").Value = 0.08939849373448538
$ws.Range("G11").Value = 0.002465057333686627
$ws.Range("I11").Value = 0.9608506379503012
$ws.Range("K11").Value = 0.612265508377277
$ws.Range("L11").Value = 0.2540730990462805
$ws.Range("N11").Value = 1.781478798012165
$ws.Range("O11").Value = 4.004736533191476

$ws.Range("B12").Value = 0.8655126568740457
$ws.Range("C12").Value = 0.08140058067841949
$ws.Range("D12").Value = 0.1156337850541576
$ws.Range("E12").Value = 0.08965291044611234
$ws.Range("G12").Value = 0.002464396977696421
$ws.Range("I12").Value = 0.9602459247985067
$ws.Range("K12").Value = 0.6248530623652471
$ws.Range("L12").Value = 0.2566326758841342
$ws.Range("N12").Value = 1.776806822199227
$ws.Range("O12").Value = 4.005757108476786

$ws.Range("B13").Value = 0.8625898826962555
$ws.Range("C13").Value = 0.08123036274531614
$ws.Range("D13").Value = 0.1151297780595826
$ws.Range("E13").Value = 0.08959784077229571
$ws.Range("G13").Value = 0.002464538620706962
$ws.Range("I13").Value = 0.9603731705083902
$ws.Range("K13").Value = 0.622141518132338
$ws.Range("L13").Value = 0.2560806419171797
$ws.Range("N13").Value = 1.777809088502035
$ws.Range("O13").Value = 4.005523908705527

$ws.Range("B14").Value = 0.8530638100548629
$ws.Range("C14").Value = 0.08067470068063187
$ws.Range("D14").Value = 0.1134870012421203
$ws.Range("E14").Value = 0.08941930013623889
$ws.Range("G14").Value = 0.002465002745802608
$ws.Range("I14").Value = 0.9607995404251284
$ws.Range("K14").Value = 0.6133008277543865
$ws.Range("L14").Value = 0.2542833238955353
$ws.Range("N14").Value = 1.781092662145557
$ws.Range("O14").Value = 4.004814457913511

$ws.Range("B15").Value = 0.8472332296036598
$ws.Range("C15").Value = 0.08033392348235679
$ws.Range("D15").Value = 0.1124814445620217
$ws.Range("E15").Value = 0.08931074871570033
$ws.Range("G15").Value = 0.002465288725056923
$ws.Range("I15").Value = 0.9610694598523111
$ws.Range("K15").Value = 0.607887383630981
$ws.Range("L15").Value = 0.2531847095638255
$ws.Range("N15").Value = 1.783115444127377
$ws.Range("O15").Value = 4.004419137853631

$ws.Range("B16").Value = 0.8138886684145064
$ws.Range("C16").Value = 0.07837466553844763
$ws.Range("D16").Value = 0.1067296190819462
$ws.Range("E16").Value = 0.08870111058090302
$ws.Range("G16").Value = 0.002466953488813428
$ws.Range("I16").Value = 0.9627504652322045
$ws.Range("K16").Value = 0.57689207826661
$ws.Range("L16").Value = 0.2469242236952169
$ws.Range("N16").Value = 1.794883371957681
$ws.Range("O16").Value = 4.002755156729449

$ws.Range("B17").Value = 0.793494185056403
$ws.Range("C17").Value = 0.07716701995828146
$ws.Range("D17").Value = 0.1032105455642522
$ws.Range("E17").Value = 0.08833821954153365
$ws.Range("G17").Value = 0.002467997956783877
$ws.Range("I17").Value = 0.9639014706282225
$ws.Range("K17").Value = 0.5579019237535476
$ws.Range("L17").Value = 0.2431151131262652
$ws.Range("N17").Value = 1.802259765055474
$ws.Range("O17").Value = 4.002271044191701

$ws.Range("B18").Value = 0.7817865619923055
$ws.Range("C18").Value = 0.07647032100832973
$ws.Range("D18").Value = 0.1011899590133254
$ws.Range("E18").Value = 0.08813357982455017
$ws.Range("G18").Value = 0.002468607248804104
$ws.Range("I18").Value = 0.9646075143190203
$ws.Range("K18").Value = 0.5469884409034194
$ws.Range("L18").Value = 0.2409358112072368
$ws.Range("N18").Value = 1.806560200097229
$ws.Range("O18").Value = 4.002189830286483

$ws.Range("B19").Value = 0.7778264880431323
$ws.Range("C19").Value = 0.07623407151199046
$ws.Range("D19").Value = 0.1005064233197288
$ws.Range("E19").Value = 0.08806499456900951
$ws.Range("G19").Value = 0.002468815013288331
$ws.Range("I19").Value = 0.9648541286620755
$ws.Range("K19").Value = 0.5432949103378917
$ws.Range("L19").Value = 0.2401999311696557
$ws.Range("N19").Value = 1.808026166303444
$ws.Range("O19").Value = 4.002196208367877

$ws.Range("B20").Value = 0.7956628633522485
$ws.Range("C20").Value = 0.07729579259617481
$ws.Range("D20").Value = 0.1035847957230231
$ws.Range("E20").Value = 0.08837642719284133
$ws.Range("G20").Value = 0.002467885887986375
$ws.Range("I20").Value = 0.9637743887394024
$ws.Range("K20").Value = 0.5599225146063418
$ws.Range("L20").Value = 0.2435193996574299
$ws.Range("N20").Value = 1.801468560173562
$ws.Range("O20").Value = 4.002302165830059

$ws.Range("B21").Value = 0.8558609297641908
$ws.Range("C21").Value = 0.08083799925509538
$ws.Range("D21").Value = 0.1139693805282036
$ws.Range("E21").Value = 0.08947157310294429
$ws.Range("G21").Value = 0.002464866069241575
$ws.Range("I21").Value = 0.9606724806160329
$ws.Range("K21").Value = 0.6158971923573802
$ws.Range("L21").Value = 0.2548107616522657
$ws.Range("N21").Value = 1.780125800597693
$ws.Range("O21").Value = 4.005014663453068

$ws.Range("B22").Value = 0.8954019695483453
$ws.Range("C22").Value = 0.08313432966168932
$ws.Range("D22").Value = 0.1207872748751839
$ws.Range("E22").Value = 0.09022357736069964
$ws.Range("G22").Value = 0.00246296809257998
$ws.Range("I22").Value = 0.9590370905372581
$ws.Range("K22").Value = 0.6525579265957333
$ws.Range("L22").Value = 0.2622930949889337
$ws.Range("N22").Value = 1.76669164711317
$ws.Range("O22").Value = 4.008543695876398

$ws.Range("B23").Value = 0.8742802578047986
$ws.Range("C23").Value = 0.08191044980260642
$ws.Range("D23").Value = 0.117145608900131
$ws.Range("E23").Value = 0.08981890622027677
$ws.Range("G23").Value = 0.002463974176332085
$ws.Range("I23").Value = 0.9598740746902621
$ws.Range("K23").Value = 0.6329844285704667
$ws.Range("L23").Value = 0.258290253245363
$ws.Range("N23").Value = 1.773814599120016
$ws.Range("O23").Value = 4.006499491529809

$ws.Range("B24").Value = 0.7946823495308877
$ws.Range("C24").Value = 0.07723758198470421
$ws.Range("D24").Value = 0.1034155892032942
$ws.Range("E24").Value = 0.08835914108056642
$ws.Range("G24").Value = 0.002467936526866309
$ws.Range("I24").Value = 0.9638317044192917
$ws.Range("K24").Value = 0.5590089923310018
$ws.Range("L24").Value = 0.2433365886745662
$ws.Range("N24").Value = 1.801826078100019
$ws.Range("O24").Value = 4.002287481746976

$ws.Range("B25").Value = 0.7096949591654607
$ws.Range("C25").Value = 0.07211713409517273
$ws.Range("D25").Value = 0.08873903310269782
$ws.Range("E25").Value = 0.0869408101769622
$ws.Range("G25").Value = 0.002472537378274637
$ws.Range("I25").Value = 0.969783490696706
$ws.Range("K25").Value = 0.4795661346940676
$ws.Range("L25").Value = 0.2276502684671868
$ws.Range("N25").Value = 1.834247881161907
$ws.Range("O25").Value = 4.005284731707576
